# tasks.xlsx - "functioneaza si DELETE TASK (interf+Excel)"
#
# This session:
#  1) deletes the "Task2.1 / do the view" row (row 4) — a task got removed,
#  2) deletes a goal ("anotherGoal") together with its task ("newGoal"/1|1)
#     and also removes two stray tasks under "newGoal" ("gsaga"/2|2 and
#     "hopa"/4|4) that were orphaned by the delete, replacing them with a
#     single new task "aTaskBecauseTheOtherWasDeleted" (progress 2|30),
#  3) adds a brand new goal "brandNewGoal" at the end of the list.
#
# The interface re-syncs the whole table back to Excel after every change,
# so every surviving date cell gets its number format re-applied too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) delete the "Task2.1 / do the view" row ---------------------------
$ws.Rows(4).Delete()

# --- 2) delete "gsaga", "hopa", "anotherGoal", "newGoal"(1|1) -------------
# (after the row-4 delete these four rows landed at 9..12)
$ws.Range("A9:G12").EntireRow.Delete()

# --- make room for the replacement task right above "uguigiuhiuh" --------
$ws.Rows(9).Insert()

$ws.Cells.Item(9, 1).Value = "Task"
$ws.Cells.Item(9, 2).Value = "aTaskBecauseTheOtherWasDeleted"
$ws.Cells.Item(9, 3).Value = 44246.84532611111
$ws.Cells.Item(9, 3).NumberFormat = "dd/MM/yyyy"
$ws.Cells.Item(9, 4).Value = "2|30"
$ws.Cells.Item(9, 5).Value = "0|0"
$ws.Cells.Item(9, 6).Value = "0%"
$ws.Cells.Item(9, 7).Value = "Just Started"

# --- 3) append the new goal "brandNewGoal" at the bottom ------------------
$ws.Cells.Item(12, 1).Value = "Goal"
$ws.Cells.Item(12, 2).Value = "brandNewGoal"
$ws.Cells.Item(12, 3).Value = 44238.814390023152
$ws.Cells.Item(12, 3).NumberFormat = "dd/mm/yyyy"
$ws.Cells.Item(12, 4).Value = "0|0"
$ws.Cells.Item(12, 5).Value = "0|0"
$ws.Cells.Item(12, 6).Value = "0%"
$ws.Cells.Item(12, 7).Value = "Just Started"

# --- the interface re-syncs every surviving date cell's display format ---
foreach ($r in 4, 5, 6, 7, 8, 10, 11) {
    $ws.Cells.Item($r, 3).NumberFormat = "dd/mm/yyyy"
}

$ws.Range("B10").Select()
